$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.397.64"
$ws.Range("E2").Value = "  -1.90%  "

$ws.Range("D3").Value = "1.837.05"
$ws.Range("E3").Value = "  -2.25%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "259.86"
$ws.Range("E5").Value = "  -6.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5192"
$ws.Range("E7").Value = "  -1.53%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3250"
$ws.Range("E8").Value = "  -6.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06759"
$ws.Range("E9").Value = "  -2.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.65"
$ws.Range("E10").Value = "  -7.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7637"
$ws.Range("E11").Value = "  -5.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07657"
$ws.Range("E12").Value = "  -2.81%  "

$ws.Range("D13").Value = "1.874.98"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.69"
$ws.Range("E14").Value = "  -1.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.019"
$ws.Range("E15").Value = "  -2.90%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.98"
$ws.Range("E17").Value = "  -4.11%  "

$ws.Range("E18").Value = "  +0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007881"
$ws.Range("E19").Value = "  -2.62%  "

$ws.Range("D20").Value = "26.430.74"
$ws.Range("E20").Value = "  -1.96%  "

$ws.Range("D21").Value = "2.067.78"
$ws.Range("E21").Value = "  -2.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.556"
$ws.Range("E22").Value = "  -4.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.436"
$ws.Range("E23").Value = "  -6.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.938"
$ws.Range("E24").Value = "  -4.14%  "

# Row 25 content swapped/changed
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.24"
$ws.Range("E25").Value = "  -2.02%  "

# Row 26 content swapped/changed
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.242"
$ws.Range("E26").Value = "  -4.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.642"
$ws.Range("E27").Value = "  -1.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.94"
$ws.Range("E28").Value = "  -3.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.38"
$ws.Range("E29").Value = "  -2.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.178"
$ws.Range("E30").Value = "  -4.47%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.143"
$ws.Range("E31").Value = "  -4.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08699"
$ws.Range("E32").Value = "  -2.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04788"
$ws.Range("E33").Value = "  -3.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.122"
$ws.Range("E34").Value = "  -4.96%  "

$ws.Range("E35").Value = "  -1.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6916"
$ws.Range("E36").Value = "  -6.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.059"
$ws.Range("E37").Value = "  -7.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01761"
$ws.Range("E38").Value = "  -5.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.193"
$ws.Range("E39").Value = "  -8.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4828"
$ws.Range("E40").Value = "  -6.44%  "

# Row 41 content swapped/changed
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "110.77"
$ws.Range("E41").Value = "  -4.81%  "

# Row 42 content swapped/changed
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8887"
$ws.Range("E42").Value = "  -7.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.093"
$ws.Range("E43").Value = "  -2.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.689"
$ws.Range("E45").Value = "  -5.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4139"
$ws.Range("E46").Value = "  -8.74%  "

# Row 47 content swapped/changed
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05858"
$ws.Range("E47").Value = "  -1.50%  "

# Row 48 content swapped/changed
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.998"
$ws.Range("E48").Value = "  -4.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1237"
$ws.Range("E49").Value = "  -8.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.77"
$ws.Range("E50").Value = "  -4.68%  "

$ws.Range("E51").Value = "  -0.73%  "
